$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.927.86'
$ws.Range('E2').Value = '  +0.59%  '

$ws.Range('D3').Value = '1.811.31'
$ws.Range('E3').Value = '  +1.68%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.25%  '

$ws.Range('E6').Value = '  -0.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4975'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.72%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3880'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.99%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09624'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +23.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.101'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.49%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.97'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.458'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.16%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.003'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.08%  '

$ws.Range('D15').Value = '1.811.50'
$ws.Range('E15').Value = '  +1.83%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.293'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.83%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.20%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.70%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06633'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.37%  '

$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.94%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.922'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('D23').Value = '27.982.21'
$ws.Range('E23').Value = '  +0.58%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.45%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.247'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.27%  '

$ws.Range('D27').Value = '2.022.68'
$ws.Range('E27').Value = '  +1.91%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.392'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.98%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.71%  '

$ws.Range('E31').Value = '  -1.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.039'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.570'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.76%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.628'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06707'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.69%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.951'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.04%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02328'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.27%  '

$ws.Range('E38').Value = '  +1.10%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.940'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.34%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6184'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.55%  '

$ws.Range('E42').Value = '  -0.17%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.146'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.43%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5872'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.32%  '

$ws.Range('E46').Value = '  -0.57%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.275'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.21%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.933'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.13%  '

$ws.Range('E50').Value = '  -2.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06791'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.22%  '
